$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    ,@(' Abu Dhabi',' September 23 2020','Mumbai won by 49 runs','Kolkata Knight Riders','Mumbai Indians','Dinesh Karthik' + [char]0x00A0 + '(c)†','30','23','5','0','130.43')
    ,@(' Abu Dhabi',' September 26 2020','KKR won by 7 wickets (with 12 balls remaining)','Kolkata Knight Riders','Sunrisers Hyderabad','Dinesh Karthik' + [char]0x00A0 + '(c)†','0','3','0','0','0.00')
    ,@(' Dubai (DSC)',' September 30 2020','KKR won by 37 runs','Kolkata Knight Riders','Rajasthan Royals','Dinesh Karthik' + [char]0x00A0 + '(c)†','1','3','0','0','33.33')
    ,@(' Abu Dhabi',' October 10 2020','KKR won by 2 runs','Kolkata Knight Riders','Kings XI Punjab','Dinesh Karthik' + [char]0x00A0 + '(c)†','58','29','8','2','200.00')
    ,@(' Sharjah',' October 03 2020','Capitals won by 18 runs','Kolkata Knight Riders','Delhi Capitals','Dinesh Karthik' + [char]0x00A0 + '(c)†','6','8','0','0','75.00')
    ,@(' Sharjah',' October 12 2020','RCB won by 82 runs','Kolkata Knight Riders','Royal Challengers Bangalore','Dinesh Karthik' + [char]0x00A0 + '(c)†','1','2','0','0','50.00')
    ,@(' Abu Dhabi',' October 07 2020','KKR won by 10 runs','Kolkata Knight Riders','Chennai Super Kings','Dinesh Karthik' + [char]0x00A0 + '(c)†','12','11','1','0','109.09')
)

$startRow = 9
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $data[$i]
    $r = $startRow + $i
    for ($c = 1; $c -le $row.Count; $c++) {
        $cell = $ws.Cells.Item($r, $c)
        if ($c -ge 7) {
            $cell.NumberFormat = "@"
        }
        $cell.Value = $row[$c - 1]
    }
}
